# Update the worksheet date in the title paragraph.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-12-29 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-12-30 Saturday", 2)

# Update each division expression in the practice table. Cells are addressed
# directly by row/column (rather than a global Find/Replace) because several
# of the old/new values repeat across the table (e.g. "69÷6=" and "27÷9="
# are both a source and a target elsewhere), which would risk cascading
# replacements with a single Replace-All pass.
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text  = "15÷8="
$t.Rows.Item(1).Cells.Item(2).Range.Text  = "50÷5="
$t.Rows.Item(1).Cells.Item(3).Range.Text  = "62÷4="
$t.Rows.Item(1).Cells.Item(4).Range.Text  = "84÷3="
$t.Rows.Item(1).Cells.Item(5).Range.Text  = "52÷7="

$t.Rows.Item(5).Cells.Item(1).Range.Text  = "69÷6="
$t.Rows.Item(5).Cells.Item(2).Range.Text  = "83÷7="
$t.Rows.Item(5).Cells.Item(3).Range.Text  = "96÷6="
$t.Rows.Item(5).Cells.Item(4).Range.Text  = "89÷2="
$t.Rows.Item(5).Cells.Item(5).Range.Text  = "97÷4="

$t.Rows.Item(9).Cells.Item(1).Range.Text  = "11÷3="
$t.Rows.Item(9).Cells.Item(2).Range.Text  = "87÷5="
$t.Rows.Item(9).Cells.Item(3).Range.Text  = "21÷6="
$t.Rows.Item(9).Cells.Item(4).Range.Text  = "40÷5="
$t.Rows.Item(9).Cells.Item(5).Range.Text  = "27÷9="

$t.Rows.Item(13).Cells.Item(1).Range.Text = "47÷4="
$t.Rows.Item(13).Cells.Item(2).Range.Text = "72÷9="
$t.Rows.Item(13).Cells.Item(3).Range.Text = "50÷7="
$t.Rows.Item(13).Cells.Item(4).Range.Text = "94÷4="
$t.Rows.Item(13).Cells.Item(5).Range.Text = "50÷3="

$t.Rows.Item(17).Cells.Item(1).Range.Text = "26÷5="
$t.Rows.Item(17).Cells.Item(2).Range.Text = "76÷9="
$t.Rows.Item(17).Cells.Item(3).Range.Text = "26÷7="
$t.Rows.Item(17).Cells.Item(4).Range.Text = "76÷9="
$t.Rows.Item(17).Cells.Item(5).Range.Text = "44÷7="
